# SCD0024-001 - Penambahan Role
# Rename sheet, update TC_ID cell, move selection and widen column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet was renamed from SCD0322 to SCD0024
$ws.Name = "SCD0024"

# TC_ID value updated from DGS-337 to SCD0024-001
$ws.Range("B2").Value = "SCD0024-001"

# Column B widened to fit the longer TC_ID text
$ws.Columns.Item(2).ColumnWidth = 11.666666666666666

# Active cell / selection moved to B3
$ws.Range("B3").Select()
